$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = 47290
$ws.Range("B2").Value = "Anthony Silva"
$ws.Range("C2").Value = "TI"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45088
$ws.Range("G2").Value = 10445.96

# Row 3
$ws.Range("A3").Value = 78551
$ws.Range("B3").Value = "Ana Carolina Lopes"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45089
$ws.Range("G3").Value = 4324.66

# Row 4
$ws.Range("A4").Value = 24814
$ws.Range("B4").Value = "Larissa Barros"
$ws.Range("C4").Value = "Vendas"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45091
$ws.Range("G4").Value = 3918.05

# Row 5
$ws.Range("A5").Value = 58971
$ws.Range("B5").Value = "Juliana Freitas"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45103
$ws.Range("G5").Value = 4597.36

# Row 6
$ws.Range("A6").Value = 69678
$ws.Range("B6").Value = "Sophie Pinto"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 45098
$ws.Range("G6").Value = 2683.2

# Row 7
$ws.Range("A7").Value = 8263
$ws.Range("B7").Value = "Arthur Novaes"
$ws.Range("C7").Value = "Operações"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45083
$ws.Range("G7").Value = 2971.39

# Row 8
$ws.Range("A8").Value = 36935
$ws.Range("B8").Value = "Sarah Ferreira"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45080
$ws.Range("G8").Value = 4389.81

# Row 9
$ws.Range("A9").Value = 7096
$ws.Range("B9").Value = "Natália Costela"
$ws.Range("C9").Value = "Jurídico"
$ws.Range("D9").Value = "Viagem de negócios"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45093
$ws.Range("G9").Value = 5236.63

# Row 10
$ws.Range("A10").Value = 74301
$ws.Range("B10").Value = "Camila Rezende"
$ws.Range("C10").Value = "Jurídico"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45086
$ws.Range("G10").Value = 7927.06

# Row 11
$ws.Range("A11").Value = 18749
$ws.Range("B11").Value = "Benjamin Barbosa"
$ws.Range("C11").Value = "Marketing"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45092
$ws.Range("G11").Value = 7196.21
